$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "correct" (right) marks row and recompute the total
# correct / total-marks summary on the concise marksheet.
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 105
$ws.Range("E12").Value = "105/140"
